# Update countries & provincias Spain
# Refresh the COVID "Pais" dashboard: bump the "last updated" timestamp,
# refresh several countries' case counts, and re-rank three country pairs
# whose totals crossed over (Nigeria/Rumania, Uruguay/Georgia, Togo/Crucero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---------------------------------------------------
$ws.Range('A1').Value = 'Datos actualizados a 13 de Julio de 2020 a las 01:50'

# --- Estados Unidos (row 4) ---------------------------------------------
$ws.Range('B4').Value = 3413060
$ws.Range('C4').Value = 57414
$ws.Range('D4').Value = 1516031
$ws.Range('E4').Value = 1759247
$ws.Range('G4').Value = 380
$ws.Range('H4').Value = 137782

# --- Canada (row 23) ------------------------------------------------------
$ws.Range('B23').Value = 107590
$ws.Range('C23').Value = 243
$ws.Range('E23').Value = 27340

# --- Rumania / Nigeria swap positions (rows 51-52) -----------------------
# Nigeria's total now edges ahead of Rumania's, so the two rows trade
# country names along with the rest of their figures.
$ws.Range('A51').Value = 'Nigeria'
$ws.Range('B51').Value = 32558
$ws.Range('C51').Value = 571
$ws.Range('D51').Value = 13447
$ws.Range('E51').Value = 18371
$ws.Range('G51').Value = 16
$ws.Range('H51').Value = 740

$ws.Range('A52').Value = 'Rumania'
$ws.Range('B52').Value = 32535
$ws.Range('C52').Value = 456
$ws.Range('D52').Value = 21545
$ws.Range('E52').Value = 9106
$ws.Range('G52').Value = 13
$ws.Range('H52').Value = 1884

# --- Japon (row 59) --------------------------------------------------------
$ws.Range('B59').Value = 21502
$ws.Range('C59').Value = 373
$ws.Range('D59').Value = 18003
$ws.Range('E59').Value = 2517

# --- Chequia (row 68) -------------------------------------------------------
$ws.Range('B68').Value = 13174
$ws.Range('C68').Value = 59
$ws.Range('D68').Value = 8247
$ws.Range('E68').Value = 4575

# --- Noruega (row 78) --------------------------------------------------------
$ws.Range('B78').Value = 8981
$ws.Range('C78').Value = 4
$ws.Range('E78').Value = 591

# --- Mauritania (row 94) -----------------------------------------------------
$ws.Range('B94').Value = 5355
$ws.Range('C94').Value = 80
$ws.Range('D94').Value = 2363
$ws.Range('E94').Value = 2845

# --- Zambia (row 120) ---------------------------------------------------------
$ws.Range('D120').Value = 1412
$ws.Range('E120').Value = 441

# --- Burkina Faso (row 139) ----------------------------------------------------
$ws.Range('B139').Value = 1036
$ws.Range('C139').Value = 3
$ws.Range('E139').Value = 114

# --- Georgia / Uruguay swap positions (rows 143-144) ---------------------------
$ws.Range('A143').Value = 'Uruguay'
$ws.Range('B143').Value = 987
$ws.Range('C143').Value = 1
$ws.Range('D143').Value = 896
$ws.Range('E143').Value = 60
$ws.Range('G143').Value = 1
$ws.Range('H143').Value = 31

$ws.Range('A144').Value = 'Georgia'
$ws.Range('C144').Value = 5
$ws.Range('D144').Value = 857
$ws.Range('E144').Value = 114
$ws.Range('H144').Value = 15

# --- Zimbabue (row 145) ----------------------------------------------------------
$ws.Range('B145').Value = 985
$ws.Range('C145').Value = 3
$ws.Range('D145').Value = 328
$ws.Range('E145').Value = 639

# --- Santo Tome y Principe (row 151) ----------------------------------------------
$ws.Range('B151').Value = 729
$ws.Range('C151').Value = 2
$ws.Range('D151').Value = 286

# --- Crucero / Togo swap positions (rows 152-153) ---------------------------------
$ws.Range('A152').Value = 'Togo'
$ws.Range('B152').Value = 720
$ws.Range('C152').Value = 10
$ws.Range('D152').Value = 513
$ws.Range('E152').Value = 192
$ws.Range('H152').Value = 15

$ws.Range('A153').Value = 'Crucero'
$ws.Range('B153').Value = 712
$ws.Range('D153').Value = 651
$ws.Range('E153').Value = 48
$ws.Range('H153').Value = 13

# --- Guyana (row 167) -----------------------------------------------------------------
$ws.Range('B167').Value = 297
$ws.Range('C167').Value = 6
$ws.Range('D167').Value = 154
